$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the obsolete earliest-date column (old column B) - shifts all header dates and data left by one column
$ws.Columns("B:B").Delete()

# Step 2: remove the two trailing forecast-origin rows that are no longer produced (old rows 23 and 24)
$ws.Rows("23:24").Delete()

# Step 3: clear stale leading forecast values that the bugfix drops from each row
$cellsToClear = @(
  "B3", "C3", "D3", "B4", "C4", "D4", "E4", "F4", "B5", "C5", "D5", "E5", "F5", "G5", "H5", "D6", "E6", "F6", "G6", "H6", "I6", "J6", "F7", "G7", "H7", "I7", "J7", "H8", "I8", "J8", "J9", "K9", "L10", "M10", "N11", "O11", "P12", "Q12", "Q13", "R13", "S13", "S14", "T14", "U14", "V14", "U15", "V15", "W15", "X15", "Y15", "Z15", "X16", "Y16", "Z16", "AA16", "AB16", "AC16", "AD16", "AB17", "AC17", "AD17", "AE17", "AF17", "AG17", "AF18", "AG18", "AH18", "AI18", "AJ18", "AK18", "AJ19", "AK19", "AL19", "AM19", "AN19", "AO19", "AN20", "AO20", "AP20", "AQ20", "AR20", "AS20", "AR21", "AS21", "AT21", "AU21", "AV21", "AW21", "AV22", "AW22", "AX22", "AY22", "AZ22"
)
foreach ($addr in $cellsToClear) {
  $ws.Range($addr).ClearContents()
}

# Step 4: recomputed forecast values for the affected cells
$ws.Range("E3").Value = 0.1715429114845124
$ws.Range("F3").Value = 0.1715429114845124
$ws.Range("G3").Value = 0.1715429114845124
$ws.Range("H3").Value = 0.1715429114845124
$ws.Range("I3").Value = 0.1715429114845124
$ws.Range("J3").Value = 0.1715429114845124
$ws.Range("K7").Value = 0.5784444854042281
$ws.Range("L7").Value = 1.133560223479058
$ws.Range("K8").Value = 1.962049292219414
$ws.Range("L8").Value = 1.985690391709771
$ws.Range("M8").Value = 2.529895848567842
$ws.Range("N8").Value = 3.633318781899142
$ws.Range("L9").Value = 2.090102686531425
$ws.Range("M9").Value = 2.205381251914007
$ws.Range("N9").Value = 2.715291551682419
$ws.Range("O9").Value = 4.060884847379076
$ws.Range("P9").Value = 3.057638025163611
$ws.Range("N10").Value = 2.448864397591044
$ws.Range("O10").Value = 2.792143403677905
$ws.Range("P10").Value = 2.42782168586293
$ws.Range("Q10").Value = 2.270469368501771
$ws.Range("R10").Value = 2.319057151538662
$ws.Range("P11").Value = 2.543955481275106
$ws.Range("Q11").Value = 2.507859322024841
$ws.Range("R11").Value = 2.508920621023392
$ws.Range("S11").Value = 2.467161166346266
$ws.Range("T11").Value = 2.536029549059826
$ws.Range("R12").Value = 2.526834392238175
$ws.Range("S12").Value = 2.518575433256176
$ws.Range("T12").Value = 2.546671316138061
$ws.Range("U12").Value = 2.480855794925163
$ws.Range("V12").Value = 3.025024236774643
$ws.Range("W12").Value = 3.120740332206995
$ws.Range("X12").Value = 3.279355759764568
$ws.Range("T13").Value = 2.530440776250154
$ws.Range("U13").Value = 2.518755579319643
$ws.Range("V13").Value = 2.69389938681992
$ws.Range("W13").Value = 2.775533179497169
$ws.Range("X13").Value = 3.107596903291299
$ws.Range("Y13").Value = 3.221757900820066
$ws.Range("Z13").Value = 2.945303709067959
$ws.Range("AA13").Value = 2.891533899000343
$ws.Range("AB13").Value = 2.827707622797226
$ws.Range("W14").Value = 2.604201945499174
$ws.Range("X14").Value = 2.657071530429667
$ws.Range("Y14").Value = 2.678174398932609
$ws.Range("Z14").Value = 2.591074440292807
$ws.Range("AA14").Value = 2.545843589346886
$ws.Range("AB14").Value = 2.413544192054795
$ws.Range("AC14").Value = 2.631992339577627
$ws.Range("AD14").Value = 2.552688975800033
$ws.Range("AE14").Value = 2.618329006605924
$ws.Range("AA15").Value = 2.623024301937549
$ws.Range("AB15").Value = 2.60322048149817
$ws.Range("AC15").Value = 2.640819364776803
$ws.Range("AD15").Value = 2.616345720823721
$ws.Range("AE15").Value = 2.671430903007876
$ws.Range("AF15").Value = 1.691013991470625
$ws.Range("AG15").Value = 2.153309886824961
$ws.Range("AH15").Value = 2.130407351599706
$ws.Range("AI15").Value = 2.137626121054947
$ws.Range("AE16").Value = 2.630644791314363
$ws.Range("AF16").Value = 2.558570068847144
$ws.Range("AG16").Value = 2.761341020331276
$ws.Range("AH16").Value = 2.785334366326175
$ws.Range("AI16").Value = 2.891950990452763
$ws.Range("AJ16").Value = 3.941556826710224
$ws.Range("AK16").Value = 4.667362054855917
$ws.Range("AL16").Value = 5.037171918133976
$ws.Range("AM16").Value = 4.951039758187648
$ws.Range("AH17").Value = 2.764442819703916
$ws.Range("AI17").Value = 2.919819837356252
$ws.Range("AJ17").Value = 3.134394395265594
$ws.Range("AK17").Value = 3.327089769540992
$ws.Range("AL17").Value = 3.641364543513781
$ws.Range("AM17").Value = 3.481452844954491
$ws.Range("AN17").Value = 2.845322256798233
$ws.Range("AO17").Value = 3.305715257492858
$ws.Range("AP17").Value = 3.153537734543965
$ws.Range("AQ17").Value = 2.838865660558509
$ws.Range("AL18").Value = 3.226452504784616
$ws.Range("AM18").Value = 3.20116940334636
$ws.Range("AN18").Value = 3.029622899744266
$ws.Range("AO18").Value = 3.019047171689593
$ws.Range("AP18").Value = 2.935215611250452
$ws.Range("AQ18").Value = 2.377254777217375
$ws.Range("AR18").Value = 2.138412043368865
$ws.Range("AS18").Value = 1.757655717321982
$ws.Range("AT18").Value = 1.831762447564067
$ws.Range("AU18").Value = 1.625773169906108
$ws.Range("AP19").Value = 2.986397903652205
$ws.Range("AQ19").Value = 2.82910658530624
$ws.Range("AR19").Value = 2.714478023861111
$ws.Range("AS19").Value = 2.633539027099796
$ws.Range("AT19").Value = 2.69124964061378
$ws.Range("AU19").Value = 2.42082970885531
$ws.Range("AV19").Value = 2.01742511619909
$ws.Range("AW19").Value = 2.159361127638926
$ws.Range("AX19").Value = 2.104676416355189
$ws.Range("AY19").Value = 2.030491763452114
$ws.Range("AT20").Value = 2.711808184127418
$ws.Range("AU20").Value = 2.79751891585911
$ws.Range("AV20").Value = 2.787508609954714
$ws.Range("AW20").Value = 2.780289798993185
$ws.Range("AX20").Value = 2.754798876280251
$ws.Range("AY20").Value = 2.559374235215039
$ws.Range("AZ20").Value = 2.733459627814305
$ws.Range("AX21").Value = 2.825169002342753
$ws.Range("AY21").Value = 2.70120649680623
$ws.Range("AZ21").Value = 2.76671919604734
